$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Status column text updates (applies everywhere "Ready for handoff" appeared) ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# --- zh-cn: Latest Target File / Latest Handback File / Latest Handback DateTime ---
$zhcn.Range("I2").Value = "1d789ca0-6a6f-4b96-ad9c-f5dd315b80d9.md"
$zhcn.Range("J2").Value = "1d789ca0-6a6f-4b96-ad9c-f5dd315b80d9.d45cbbe64a3495b53f9fb62a53b351d25430b9c4.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-10-20 01:40:11"

$zhcn.Range("I3").Value = "c8d9c541-ee65-49cf-bba5-f5cb9c7f4535.md"
$zhcn.Range("J3").Value = "c8d9c541-ee65-49cf-bba5-f5cb9c7f4535.3f7a505a287b0caceda910459b4fc35232b83c9a.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-10-20 01:40:11"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d96f5ce070753ae6d58b155b85a04eb519c52b0c/e2e/1d789ca0-6a6f-4b96-ad9c-f5dd315b80d9.md", "", "", "1d789ca0-6a6f-4b96-ad9c-f5dd315b80d9.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d96f5ce070753ae6d58b155b85a04eb519c52b0c/e2e/c8d9c541-ee65-49cf-bba5-f5cb9c7f4535.md", "", "", "c8d9c541-ee65-49cf-bba5-f5cb9c7f4535.md") | Out-Null

$zhcn.Range("I2:I3").Style = "Hyperlink"

# --- de-de: Latest Target File / Latest Handback File / Latest Handback DateTime ---
$dede.Range("I2").Value = "1d789ca0-6a6f-4b96-ad9c-f5dd315b80d9.md"
$dede.Range("J2").Value = "1d789ca0-6a6f-4b96-ad9c-f5dd315b80d9.d45cbbe64a3495b53f9fb62a53b351d25430b9c4.de-de.xlf"
$dede.Range("K2").Value = "2016-10-20 01:40:29"

$dede.Range("I3").Value = "c8d9c541-ee65-49cf-bba5-f5cb9c7f4535.md"
$dede.Range("J3").Value = "c8d9c541-ee65-49cf-bba5-f5cb9c7f4535.3f7a505a287b0caceda910459b4fc35232b83c9a.de-de.xlf"
$dede.Range("K3").Value = "2016-10-20 01:40:29"

$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d96f5ce070753ae6d58b155b85a04eb519c52b0c/e2e/1d789ca0-6a6f-4b96-ad9c-f5dd315b80d9.md", "", "", "1d789ca0-6a6f-4b96-ad9c-f5dd315b80d9.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d96f5ce070753ae6d58b155b85a04eb519c52b0c/e2e/c8d9c541-ee65-49cf-bba5-f5cb9c7f4535.md", "", "", "c8d9c541-ee65-49cf-bba5-f5cb9c7f4535.md") | Out-Null

$dede.Range("I2:I3").Style = "Hyperlink"

# --- Column widths (match autofit results capturing the new, wider content) ---
$overview.Columns.Item(5).ColumnWidth = 29.1443716684978
$overview.Columns.Item(6).ColumnWidth = 29.1443716684978

$zhcn.Columns.Item(3).ColumnWidth = 29.1443716684978
$zhcn.Columns.Item(9).ColumnWidth = 39.1666666666667
$zhcn.Columns.Item(10).ColumnWidth = 39.1666666666667

$dede.Columns.Item(3).ColumnWidth = 29.1443716684978
$dede.Columns.Item(9).ColumnWidth = 39.1666666666667
$dede.Columns.Item(10).ColumnWidth = 39.1666666666667
